$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.72"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'37.66"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'1.76%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.159"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'1.28%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.07888"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'2.18%"
$ws.Range("E5").ClearFormats()
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.413"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'1.18%"
$ws.Range("E6").ClearFormats()
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.299"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'1.21%"
$ws.Range("E7").ClearFormats()
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.908"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'1.67%"
$ws.Range("E8").ClearFormats()
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.968"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-7.07%"
$ws.Range("E9").ClearFormats()
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9227"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'0.48%"
$ws.Range("E10").ClearFormats()
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1225"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'1.34%"
$ws.Range("E11").ClearFormats()
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1923"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'1.87%"
$ws.Range("E12").ClearFormats()
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09182"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'5.41%"
$ws.Range("E13").ClearFormats()
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03346"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-1.12%"
$ws.Range("E14").ClearFormats()
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09582"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-1.20%"
$ws.Range("E15").ClearFormats()
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001382"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'1.21%"
$ws.Range("E16").ClearFormats()
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005745"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'-6.13%"
$ws.Range("E17").ClearFormats()
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.519"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-1.23%"
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'0.3443"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'2.07%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'5.267"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'4.91%"
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'-0.32%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.2592"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'3.86%"
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'-0.38%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.04369"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'0.83%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.001249"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'2.81%"
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'0.004705"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'5.30%"
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'-9.71%"
$ws.Range("E27").ClearFormats()
$ws.Range("D39").Value = "'0.02315"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'4.19%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.05086"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'3.36%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007459"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-1.90%"
$ws.Range("E41").ClearFormats()
$ws.Range("B42").Value = "Dexo"
$ws.Range("C42").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D42").Value = "'0.008884"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-9.68%"
$ws.Range("E42").ClearFormats()
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1357"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'2.12%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.001922"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-3.82%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.008622"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-2.15%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00006602"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-3.21%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'-0.22%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.003351"
$ws.Range("D48").ClearFormats()
$ws.Range("E49").Value = "'-7.91%"
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'-0.22%"
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'-0.22%"
$ws.Range("E51").ClearFormats()
